$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run RU (Russia, column L) 1001; without crop
# Update the affected numeric values in rows 2-4 (B, I, L columns)

$ws.Range("B2").Value = 0.680815013747804
$ws.Range("I2").Value = 0.68337742404715
$ws.Range("L2").Value = 0.689906588508866

$ws.Range("B3").Value = 0.698563864639442
$ws.Range("L3").Value = 0.733075012505842

$ws.Range("B4").Value = 0.66248064669155
$ws.Range("L4").Value = 0.644869545559419
